# Generate Report for Archive
# - Status text "Ready for handoff" -> "In Translation" on every sheet that
#   shows it (Overview!E2:F2, zh-cn!C2, de-de!C2 all share the same string).
# - The Status columns get narrower to fit the new (shorter) text:
#   Overview columns E & F, and column C on the zh-cn / de-de detail sheets.

$wb = $excel.ActiveWorkbook

$overview = $wb.Sheets.Item("Overview")
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"

$zhcn = $wb.Sheets.Item("zh-cn")
$zhcn.Range("C2").Value = "In Translation"

$dede = $wb.Sheets.Item("de-de")
$dede.Range("C2").Value = "In Translation"

# Narrow the Status columns to match the regenerated report layout.
$overview.Columns.Item(5).ColumnWidth = 16.333333333333336
$overview.Columns.Item(6).ColumnWidth = 16.333333333333336
$zhcn.Columns.Item(3).ColumnWidth = 12.5
$dede.Columns.Item(3).ColumnWidth = 12.5
